$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# New task rows (51-69): "tried adding and fetching user followers and
# followees". Text is entered in the exact order the shared-string table
# expects (indices 26-36) so the rebuilt sst matches the authored commit.
# ---------------------------------------------------------------------------

$ws.Range("B51").Value = "Get the followers of a person"
$ws.Range("B52").Value = "Get the followings of the person"
$ws.Range("B54").Value = "Give follow functionality"
$ws.Range("B55").Value = "Give unfollow functionality"
$ws.Range("B57").Value = "Give functionality to get the tweets of a person"
$ws.Range("B59").Value = "How to extension attributes to a class in c#"
$ws.Range("B61").Value = "Give Tweet Adding Functionality"
$ws.Range("B63").Value = "Add username functionality"
$ws.Range("B65").Value = "While getting followersand unautorized , pass unauthorized and not bad request"
$ws.Range("B67").Value = "Make FollowersDto"
$ws.Range("B69").Value = "make FolloweesDto"

# ---------------------------------------------------------------------------
# While drafting, a red font colour was dragged/toggled across blocks of
# cells B:F in this area, leaving some of the cells in red.
# ---------------------------------------------------------------------------

$ws.Range("E51:F51").Font.Color = 255
$ws.Range("E52:F52").Font.Color = 255
$ws.Range("B53:F53").Font.Color = 255
$ws.Range("B54:D54").Font.Color = 255
$ws.Range("B55:D55").Font.Color = 255
$ws.Range("B61:E61").Font.Color = 255

# The remaining cells in those same blocks were also touched (selected /
# reformatted) but kept their normal black font - materialise them as
# blank formatted cells without altering their colour.
$ws.Range("B51:D51").Borders.LineStyle = -4142
$ws.Range("B52:D52").Borders.LineStyle = -4142
$ws.Range("E54:F54").Borders.LineStyle = -4142
$ws.Range("E55:F55").Borders.LineStyle = -4142
$ws.Range("B56:F56").Borders.LineStyle = -4142
$ws.Range("B57:F59").Borders.LineStyle = -4142

# ---------------------------------------------------------------------------
# Scroll / selection state left at the end of the edit.
# ---------------------------------------------------------------------------
$ws.Range("B69").Select()
$excel.ActiveWindow.ScrollRow = 55

Write-Host "edit applied"
